$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3200
$ws.Range("J9").Value = 3200
$ws.Range("L9").Value = 3200
$ws.Range("N9").Value = -3538

$ws.Range("H17").Value = 6835.5
$ws.Range("J17").Value = 7116.316
$ws.Range("L17").Value = 21348.948
$ws.Range("N17").Value = -21684.948

$ws.Range("H70").Value = 139837.62
$ws.Range("I70").Value = 342666.34
$ws.Range("K70").Value = 1027999.02
$ws.Range("M70").Value = -1027729.02

$ws.Range("H73").Value = 139837.62
$ws.Range("I73").Value = 342666.34
$ws.Range("K73").Value = 1027999.02
$ws.Range("M73").Value = -1027063.02

$ws.Range("H135").Value = 1112.3529
$ws.Range("I135").Value = 532
$ws.Range("J135").Value = 2998.5
$ws.Range("K135").Value = 4788
$ws.Range("L135").Value = 26986.5
$ws.Range("M135").Value = -2253
$ws.Range("N135").Value = -32056.5

$ws.Range("H137").Value = 12823139
$ws.Range("I137").Value = 47620716
$ws.Range("J137").Value = 2979.4036
$ws.Range("K137").Value = 142862148
$ws.Range("L137").Value = 8938.210800000001
$ws.Range("M137").Value = -142859598
$ws.Range("N137").Value = -14038.2108

$ws.Range("H138").Value = 2800.413
$ws.Range("J138").Value = 3161.2144
$ws.Range("L138").Value = 9483.643199999999
$ws.Range("N138").Value = -19763.6432

$ws.Range("H141").Value = 2549.7334
$ws.Range("I141").Value = 2549.7334
$ws.Range("K141").Value = 7649.2002
$ws.Range("M141").Value = -2469.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 90910824
$ws.Range("I45").Value = 90910824
$ws.Range("K45").Value = 90910824
$ws.Range("M45").Value = -90910447

$ws.Range("H74").Value = 10754991
$ws.Range("I74").Value = 15153194
$ws.Range("J74").Value = 3827.111
$ws.Range("K74").Value = 15153194
$ws.Range("L74").Value = 3827.111
$ws.Range("M74").Value = -15152320
$ws.Range("N74").Value = -5575.111

$ws.Range("H77").Value = 10754991
$ws.Range("I77").Value = 15153194
$ws.Range("J77").Value = 3827.111
$ws.Range("K77").Value = 75765970
$ws.Range("L77").Value = 19135.555
$ws.Range("M77").Value = -75761602
$ws.Range("N77").Value = -27871.555

$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378

$ws.Range("H132").Value = 3453.8948
$ws.Range("I132").Value = 2681.8
$ws.Range("J132").Value = 4938.6924
$ws.Range("K132").Value = 8045.400000000001
$ws.Range("L132").Value = 14816.0772
$ws.Range("M132").Value = -5515.400000000001
$ws.Range("N132").Value = -19876.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4494.033
$ws.Range("I20").Value = 4108.5713
$ws.Range("J20").Value = 5393.4443
$ws.Range("K20").Value = 4108.5713
$ws.Range("L20").Value = 5393.4443
$ws.Range("M20").Value = -3861.5713
$ws.Range("N20").Value = -5887.4443

$ws.Range("H21").Value = 69989.336
$ws.Range("J21").Value = 69989.336
$ws.Range("L21").Value = 69989.336
$ws.Range("N21").Value = -70461.336

$ws.Range("H94").Value = 2077
$ws.Range("I94").Value = 2077
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2077
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1626
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7996.25
$ws.Range("I22").Value = 6990.3335
$ws.Range("K22").Value = 6990.3335
$ws.Range("M22").Value = -6640.3335

$ws.Range("H25").Value = 1174.25
$ws.Range("I25").Value = 1174.25
$ws.Range("K25").Value = 1174.25
$ws.Range("M25").Value = -1000.25

$ws.Range("H31").Value = 28302.049
$ws.Range("I31").Value = 2120.2
$ws.Range("K31").Value = 2120.2
$ws.Range("M31").Value = -1825.2

$ws.Range("H34").Value = 28302.049
$ws.Range("I34").Value = 2120.2
$ws.Range("K34").Value = 2120.2
$ws.Range("M34").Value = -1918.2

$ws.Range("H99").Value = 2799.2
$ws.Range("J99").Value = 2799.2
$ws.Range("L99").Value = 2799.2
$ws.Range("N99").Value = -5795.2

$ws.Range("H108").Value = 82497.5
$ws.Range("J108").Value = 82497.5
$ws.Range("L108").Value = 82497.5
$ws.Range("N108").Value = -90177.5

$ws.Range("H126").Value = 2799.2
$ws.Range("J126").Value = 2799.2
$ws.Range("L126").Value = 8397.599999999999
$ws.Range("N126").Value = -13337.6

$ws.Range("H132").Value = 3613.1072
$ws.Range("I132").Value = 3091.2222
$ws.Range("K132").Value = 9273.6666
$ws.Range("M132").Value = -6743.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1428653.2
$ws.Range("I33").Value = 2500075.5
$ws.Range("J33").Value = 90.333336
$ws.Range("K33").Value = 15000453
$ws.Range("L33").Value = 542.000016
$ws.Range("M33").Value = -15000170
$ws.Range("N33").Value = -1108.000016

$ws.Range("H75").Value = 142863170
$ws.Range("J75").Value = 8082.8
$ws.Range("L75").Value = 24248.4
$ws.Range("N75").Value = -26244.4

$ws.Range("H78").Value = 142863170
$ws.Range("J78").Value = 8082.8
$ws.Range("L78").Value = 72745.2
$ws.Range("N78").Value = -82729.2

$ws.Range("H126").Value = 37039732
$ws.Range("I126").Value = 1837.25
$ws.Range("J126").Value = 66670050
$ws.Range("K126").Value = 5511.75
$ws.Range("L126").Value = 200010150
$ws.Range("M126").Value = -571.75
$ws.Range("N126").Value = -200020030

$ws.Range("H132").Value = 3525.75
$ws.Range("J132").Value = 4423.7
$ws.Range("L132").Value = 39813.3
$ws.Range("N132").Value = -44873.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 29498
$ws.Range("I21").Value = 29498
$ws.Range("K21").Value = 29498
$ws.Range("M21").Value = -29325

$ws.Range("H30").Value = 29498
$ws.Range("I30").Value = 29498
$ws.Range("K30").Value = 29498
$ws.Range("M30").Value = -29393

$ws.Range("H32").Value = 32499.75
$ws.Range("J32").Value = 32499.75
$ws.Range("L32").Value = 32499.75
$ws.Range("N32").Value = -33091.75

$ws.Range("H42").Value = 38000
$ws.Range("J42").Value = 38000
$ws.Range("L42").Value = 38000
$ws.Range("N42").Value = -38970

$ws.Range("H97").Value = 1923.7273
$ws.Range("I97").Value = 1676.1
$ws.Range("K97").Value = 1676.1
$ws.Range("M97").Value = -1180.1

$ws.Range("H113").Value = 2896.5908
$ws.Range("I113").Value = 1730
$ws.Range("K113").Value = 1730
$ws.Range("M113").Value = 440

$ws.Range("H114").Value = 59999.5
$ws.Range("J114").Value = 59999.5
$ws.Range("L114").Value = 59999.5
$ws.Range("N114").Value = -68677.5

$ws.Range("H115").Value = 38000
$ws.Range("J115").Value = 38000
$ws.Range("L115").Value = 38000
$ws.Range("N115").Value = -40350

$ws.Range("H122").Value = 5518.4
$ws.Range("I122").Value = 4042.5217
$ws.Range("J122").Value = 6775.6294
$ws.Range("K122").Value = 12127.5651
$ws.Range("L122").Value = 20326.8882
$ws.Range("M122").Value = -9677.5651
$ws.Range("N122").Value = -25226.8882

$ws.Range("H126").Value = 3493.4285
$ws.Range("I126").Value = 1699.8334
$ws.Range("K126").Value = 5099.5002
$ws.Range("M126").Value = -2629.5002

$ws.Range("H132").Value = 2549.818
$ws.Range("I132").Value = 2072.4707
$ws.Range("J132").Value = 4172.8
$ws.Range("K132").Value = 6217.4121
$ws.Range("L132").Value = 12518.4
$ws.Range("M132").Value = -3687.4121
$ws.Range("N132").Value = -17578.4

$ws.Range("H134").Value = 75869
$ws.Range("J134").Value = 75869
$ws.Range("L134").Value = 227607
$ws.Range("N134").Value = -232677

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7243.5
$ws.Range("I7").Value = 4978.4
$ws.Range("J7").Value = 10074.875
$ws.Range("K7").Value = 4978.4
$ws.Range("L7").Value = 10074.875
$ws.Range("M7").Value = -4866.4
$ws.Range("N7").Value = -10298.875

$ws.Range("H16").Value = 2436.8235
$ws.Range("I16").Value = 558.36365
$ws.Range("K16").Value = 558.36365
$ws.Range("M16").Value = -388.36365

$ws.Range("H40").Value = 7581.8667
$ws.Range("I40").Value = 6453.1577
$ws.Range("K40").Value = 6453.1577
$ws.Range("M40").Value = -6317.1577

$ws.Range("H61").Value = 6430.85
$ws.Range("I61").Value = 5347.3335
$ws.Range("K61").Value = 5347.3335
$ws.Range("M61").Value = -5145.3335

$ws.Range("H113").Value = 6430.85
$ws.Range("I113").Value = 5347.3335
$ws.Range("K113").Value = 5347.3335
$ws.Range("M113").Value = -3177.3335

$ws.Range("H126").Value = 7243.5
$ws.Range("I126").Value = 4978.4
$ws.Range("J126").Value = 10074.875
$ws.Range("K126").Value = 14935.2
$ws.Range("L126").Value = 30224.625
$ws.Range("M126").Value = -12465.2
$ws.Range("N126").Value = -35164.625

$ws.Range("H132").Value = 5544.875
$ws.Range("I132").Value = 5790.472
$ws.Range("K132").Value = 17371.416
$ws.Range("M132").Value = -14841.416

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20480

$ws.Range("H132").Value = 2277.0303
$ws.Range("I132").Value = 1687.4828
$ws.Range("K132").Value = 5062.4484
$ws.Range("M132").Value = -2532.4484
